$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31 (A column): the original "Não recebi meu saque" question text
# gets keyword tags appended to it.
$ws.Cells.Item(31, 1).Value = "Não recebi meu saque ,porque?atrasado, saque atrasado"

# --- Row-height corrections on existing wrapped rows (auto recalculated
# by Excel after the new content below was added / fonts touched).
$ws.Rows.Item(9).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 75
$ws.Rows.Item(14).RowHeight = 165
$ws.Rows.Item(15).RowHeight = 285
$ws.Rows.Item(18).RowHeight = 345
$ws.Rows.Item(19).RowHeight = 45
$ws.Rows.Item(20).RowHeight = 30
$ws.Rows.Item(22).RowHeight = 45
$ws.Rows.Item(29).RowHeight = 105
$ws.Rows.Item(31).RowHeight = 30
$ws.Rows.Item(32).RowHeight = 45
$ws.Rows.Item(33).RowHeight = 60

# --- New FAQ rows appended at the bottom of the sheet.
$ws.Cells.Item(34, 1).Value = "Quais são minhas chances de ganhar ?"
$ws.Cells.Item(34, 2).Value = "Suas chances são enormes e com muita sorte torcemos por você!"
$ws.Cells.Item(34, 2).WrapText = $true

$ws.Cells.Item(35, 1).Value = "Vocês tem bingo ?Bingo online"
$ws.Cells.Item(35, 2).Value = "Sim temos sim, segue o link :https://betmotion.com/br/bingo-online/bingo-cartelas-online"
$ws.Cells.Item(35, 2).WrapText = $true

$ws.Cells.Item(36, 1).Value = "qual o resultado da quarta de giros ?"
$ws.Cells.Item(36, 2).Value = "o resultado foi a,b,c e d"
$ws.Cells.Item(36, 2).WrapText = $true

$ws.Cells.Item(37, 1).Value = "qual o e-mail para contato ?"
$ws.Cells.Item(37, 2).WrapText = $true
$ws.Hyperlinks.Add($ws.Cells.Item(37, 2), "mailto:suporte@betmotion.com", "", "", "suporte@betmotion.com")

$ws.Cells.Item(38, 1).Value = "qual o whatsapp para contato ?"
$ws.Cells.Item(38, 2).Value = "11 91427-1299"
$ws.Cells.Item(38, 2).WrapText = $true

$ws.Cells.Item(39, 1).Value = "como mantenho contato ?"
$ws.Cells.Item(39, 2).Value = "whatsapp 11 91427-1299 e email suporte@betmotion.com"
$ws.Cells.Item(39, 2).WrapText = $true

# --- Selection moves to the last data row, matching the author's final
# cursor position after entering the new rows.
$ws.Range("A36").Select()
